# Corrects the "Fitness" values (column C) of the run_27 log sheet.
# Per the commit "correction in sa algorithm and 746 logs", column C
# (Fitness) values for generations are recomputed/corrected while
# columns A (Run) and B (Generation) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is (startRow, endRow, newFitnessValue) for column C.
# Rows 2-252 correspond to Generation 0-250 of Run 27.
$updates = @(
    @(2,   26,  7770),
    @(27,  68,  7581),
    @(69,  115, 7534),
    @(116, 116, 7320),
    @(117, 168, 7312),
    @(169, 252, 7310)
)

foreach ($u in $updates) {
    $startRow = $u[0]
    $endRow   = $u[1]
    $value    = $u[2]
    $rangeAddr = "C" + $startRow + ":C" + $endRow
    $ws.Range($rangeAddr).Value = $value
}
